# Weekly update: insert two new records at the top of the Lechuga data
# block (right after the existing header-less data at row 1032),
# pushing all subsequent rows down by two and appending the data for
# the new week (date serial 45041).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 1033 — this shifts
# the existing rows 1033:1074 down to 1035:1076 (and grows the used
# range / dimension to A1:R1076 automatically).
$ws.Rows("1033:1034").Insert()

# --- New row 1033: Conconina(o), week of 2023-04-25 (serial 45041) ---
$ws.Cells.Item(1033, 1).Value  = 7
$ws.Cells.Item(1033, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(1033, 3).Value  = "Ñuble"
$ws.Cells.Item(1033, 4).Value  = 45041
$ws.Cells.Item(1033, 5).Value  = 16
$ws.Cells.Item(1033, 6).Value  = 100112033
$ws.Cells.Item(1033, 7).Value  = "Lechuga"
$ws.Cells.Item(1033, 8).Value  = "Conconina(o)"
$ws.Cells.Item(1033, 9).Value  = "Primera"
$ws.Cells.Item(1033, 10).Value = 60
$ws.Cells.Item(1033, 11).Value = 6000
$ws.Cells.Item(1033, 12).Value = 6000
$ws.Cells.Item(1033, 13).Value = 6000
$ws.Cells.Item(1033, 14).Value = '$/caja 10 unidades'
$ws.Cells.Item(1033, 15).Value = "Región del Maule"
$ws.Cells.Item(1033, 16).Value = 600
$ws.Cells.Item(1033, 17).Value = 10
$ws.Cells.Item(1033, 18).Value = "Hortaliza"

# --- New row 1034: Escarola, week of 2023-04-25 (serial 45041) ---
$ws.Cells.Item(1034, 1).Value  = 7
$ws.Cells.Item(1034, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(1034, 3).Value  = "Ñuble"
$ws.Cells.Item(1034, 4).Value  = 45041
$ws.Cells.Item(1034, 5).Value  = 16
$ws.Cells.Item(1034, 6).Value  = 100112033
$ws.Cells.Item(1034, 7).Value  = "Lechuga"
$ws.Cells.Item(1034, 8).Value  = "Escarola"
$ws.Cells.Item(1034, 9).Value  = "Primera"
$ws.Cells.Item(1034, 10).Value = 60
$ws.Cells.Item(1034, 11).Value = 8000
$ws.Cells.Item(1034, 12).Value = 8000
$ws.Cells.Item(1034, 13).Value = 8000
$ws.Cells.Item(1034, 14).Value = '$/caja 15 unidades'
$ws.Cells.Item(1034, 15).Value = "Región del Maule"
$ws.Cells.Item(1034, 16).Value = 533
$ws.Cells.Item(1034, 17).Value = 15
$ws.Cells.Item(1034, 18).Value = "Hortaliza"
